# Salesforce Python code and other things.
#
# Applies the edits made to the "Sheet1" (3rd position tab, internal
# sheet2.xml) worksheet of the dow.xlsx workbook: a couple of label
# renames, a row of old data removed, and a new "PFS-Web" block of
# figures (columns F-I) added alongside the existing "401K" block
# (columns C-D). Finishes by making "Sheet1" the active/selected tab
# with I9 as the selected cell (matching the author's last on-screen
# state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- relabel existing entries ---------------------------------------
$ws.Range("D7").Value = "401K"
$ws.Range("D11").Value = "Off - 22 days"

# --- drop the old "safar" row (C12/D12) ------------------------------
$ws.Range("C12:D12").ClearContents()

# --- new "PFS-Web" figures (F5:G11) ----------------------------------
$ws.Range("F5").Value = "PFS-Web"
$ws.Range("F6").Value = 180
$ws.Range("F7").Formula = "=19*0.12"
$ws.Range("F10").Formula = "=SUM(F6:F9)"
$ws.Range("F11").Value = 12
$ws.Range("G11").Value = "22 days"
$ws.Range("F15").Formula = "=SUM(F10:F14)"

# --- new helper figures (I6:I8) ---------------------------------------
$ws.Range("I6").Formula = "=165*0.12"
$ws.Range("I7").Value = 165
$ws.Range("I8").Formula = "=SUM(I6:I7)"

# --- make Sheet1 the active tab/selection, matching the saved view ---
$ws.Activate() | Out-Null
$ws.Range("I9").Select() | Out-Null
